$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'crop athletic leggings'
    2 = 'kid basketball knee pads'
    3 = 'men''s spandex leggings'
    4 = 'hex pants'
    5 = 'medical compression pants'
    6 = 'fitness tights for men'
    7 = 'the rock mens basketball'
    8 = 'black men tights'
    9 = 'd man basketball'
    10 = 'asics compression pants'
    11 = 'knee armor knee pads'
    12 = 'adidas tights for men'
    13 = 'slide on knee pads'
    14 = 'men athletic tights'
    15 = 'blue mens compression pants'
    16 = 'sport tights for men'
    17 = 'compression basketball tights'
    18 = 'kids compression pants'
    19 = 'usa tights men'
    20 = 'razor knee pads'
    21 = 'knee pads addidas'
    22 = 'knee pads hunting'
    23 = 'knee pads leggings'
    24 = 'knee pads elbow pads youth'
    25 = 'knee pads slim'
    26 = 'basketball pants adidas'
    27 = 'mens tights navy'
    28 = 'mens adidas basketball pants'
    29 = 'compression pants men 3xl'
    30 = 'compression pants 2xu'
    31 = 'pro x knee pad'
    32 = 'puma compression pants men'
    33 = 'men''s basketball pants'
    34 = 'gray baseball pants youth girls'
    35 = 'mens basketball jacket'
    36 = 'photography knee pads'
    37 = 'men''s tights leggings'
    38 = 'baseball pants men grey'
    39 = 'venom compression pants'
    40 = 'padded knee tights'
    41 = 'ua compression pants'
    42 = 'men workout tights'
    43 = 'flag compression pants'
    44 = 'cool knee pads'
    45 = 'navy compression leggings'
    46 = 'force knee pads'
    47 = 'mens wrestling pants'
    48 = 'mens capri compression pants'
    49 = 'woman compression pants'
    50 = 'purple knee pads'
    51 = 'reebok knee pads'
    52 = 'venum compression pants men'
    53 = 'purple athletic leggings'
    54 = 'thermal compression pants'
    55 = 'addidas knee pads'
    56 = 'jordan mens tights'
    57 = 'purple compression pants men'
    58 = 'russell compression pants'
    59 = 'blue knee pads for basketball'
    60 = 'elbow knee pad'
    61 = 'tights mens'
    62 = 'super compression leggings'
    63 = 'mens leggins'
    64 = 'knee pads for teens'
    65 = 'green mens compression pants'
    66 = 'pants with padded knees'
    67 = 'compression with pads'
    68 = 'knee pads for sleeping'
    69 = 'mens winter compression pants'
    70 = 'knee pads for hvac'
    71 = 'yoga pants with knee pads'
    72 = 'black pants with knee pads'
    73 = 'kids compression knee pads'
    74 = 'elite basketball pants'
    75 = 'nike leggings mens'
    76 = 'compression pants men 3 pack'
    77 = 'compression pants baseball'
    78 = 'colored compression leggings'
    79 = 'mens compression pants with pockets'
    80 = 'mens compression pants xxl'
    81 = 'knee pads wheels'
    82 = 'ua basketball knee pads'
    83 = 'protective knee pad'
    84 = 'grey knee pads basketball'
    85 = 'orange knee pads for basketball'
    86 = 'adidas youth compression pants'
    87 = 'copper compression tights for men'
    88 = 'basketball knee pads youth boys mcdavid'
    89 = 'nike youth basketball knee pads'
    90 = 'nike kneepads'
    91 = 'knee pad strap'
    92 = 'internal knee pad'
    93 = 'basketball legings'
    94 = 'basketball knee sleves'
    95 = 'under armour baseball pants men'
    96 = 'cold gear compression leggings men'
    97 = 'mcgregor knee pads'
    98 = 'knee compression pants men'
    99 = 'compression knee pads for basketball for kids'
    100 = 'tesla compression pant'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $values[$row]
}